# cold_room_status.xlsx edit
# Commit: "Change utilization to use not_installed and make other model selection better"
#
# 1. In the "choices" sheet, the current_use choice option
#    in_store_for_allocation / In Store For Allocation / Almacenado Epsperando asignación
#    is replaced by
#    not_installed / Not Installed / No Instalado
# 2. The active worksheet / selection moves from "survey" to "choices"
#    (authoring UI state: the author ended their session with the choices sheet
#    active and a different cell selected there).

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$choices  = $wb.Worksheets.Item("choices")
$settings = $wb.Worksheets.Item("settings")

# --- 1. Update the current_use choice list row (row 15) ---
$choices.Range("B15").Value = "not_installed"
$choices.Range("C15").Value = "Not Installed"
$choices.Range("D15").Value = "No Instalado"

# --- 2. Update UI selection state ---
# survey keeps its existing selection (F5); settings keeps its existing
# selection (D6). The choices sheet becomes the active/selected sheet with
# a new active cell.
$choices.Activate()
$choices.Range("D27").Select()
